$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "66.678.63"
Set-TextValue "E2" "  +3.37%  "

# Row 3
Set-TextValue "D3" "3.435.42"
Set-TextValue "E3" "  +2.55%  "

# Row 4
Set-TextValue "E4" "  -0.08%  "

# Row 5
Set-TextValue "D5" "570.54"
Set-TextValue "E5" "  +2.15%  "

# Row 6
Set-TextValue "D6" "182.31"
Set-TextValue "E6" "  +3.69%  "

# Row 7
Set-TextValue "D7" "0.632"
Set-TextValue "E7" "  +1.88%  "

# Row 8
Set-TextValue "D8" "3.430.78"
Set-TextValue "E8" "  +2.65%  "

# Row 9
Set-TextValue "E9" "  -0.06%  "

# Row 10
Set-TextValue "E10" "  +3.83%  "

# Row 11
Set-TextValue "D11" "0.643"
Set-TextValue "E11" "  +1.44%  "

# Row 12
Set-TextValue "D12" "56.01"
Set-TextValue "E12" "  +4.03%  "

# Row 13
Set-TextValue "D13" "0.0000278"
Set-TextValue "E13" "  +1.14%  "

# Row 14
Set-TextValue "D14" "9.37"
Set-TextValue "E14" "  +3.21%  "

# Row 15
Set-TextValue "D15" "3.985.50"
Set-TextValue "E15" "  +2.37%  "

# Row 16
Set-TextValue "D16" "18.56"
Set-TextValue "E16" "  +1.84%  "

# Row 17
Set-TextValue "D17" "3.440.58"
Set-TextValue "E17" "  +2.69%  "

# Row 18
Set-TextValue "E18" "  +0.57%  "

# Row 19
Set-TextValue "D19" "66.676.72"
Set-TextValue "E19" "  +2.90%  "

# Row 20
Set-TextValue "D20" "12.05"
Set-TextValue "E20" "  +2.57%  "

# Row 21
Set-TextValue "D21" "1.01"
Set-TextValue "E21" "  +2.66%  "

# Row 22
Set-TextValue "D22" "484.22"
Set-TextValue "E22" "  +7.27%  "

# Row 23
Set-TextValue "D23" "16.30"
Set-TextValue "E23" "  +14.79%  "

# Row 24
Set-TextValue "D24" "5.01"
Set-TextValue "E24" "  +1.71%  "

# Row 25
Set-TextValue "D25" "4.21"
Set-TextValue "E25" "  +2.12%  "

# Row 26
Set-TextValue "D26" "89.14"
Set-TextValue "E26" "  +2.79%  "

# Row 27
Set-TextValue "D27" "2.96"
Set-TextValue "E27" "  +2.14%  "

# Row 28
Set-TextValue "D28" "10.95"
Set-TextValue "E28" "  +1.73%  "

# Row 29
Set-TextValue "D29" "9.07"
Set-TextValue "E29" "  +4.14%  "

# Row 30
Set-TextValue "D30" "31.33"
Set-TextValue "E30" "  +1.59%  "

# Row 31
Set-TextValue "D31" "7.19"
Set-TextValue "E31" "  +8.46%  "

# Row 32
Set-TextValue "D32" "593.24"
Set-TextValue "E32" "  +3.81%  "

# Row 33
Set-TextValue "D33" "11.66"
Set-TextValue "E33" "  +1.68%  "

# Row 34
Set-TextValue "D34" "63.09"
Set-TextValue "E34" "  +3.46%  "

# Row 35
Set-TextValue "D35" "0.111"
Set-TextValue "E35" "  +3.66%  "

# Row 36
Set-TextValue "B36" "Dai"
Set-TextValue "C36" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D36" "0.999"
Set-TextValue "E36" "  -0.13%  "

# Row 37
Set-TextValue "B37" "Kaspa"
Set-TextValue "C37" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D37" "0.148"
Set-TextValue "E37" "  +5.14%  "

# Row 38
Set-TextValue "E38" "  -1.93%  "

# Row 39
Set-TextValue "E39" "  +4.65%  "

# Row 40
Set-TextValue "B40" "PEPE"
Set-TextValue "C40" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D40" "0.0₃0771"
Set-TextValue "E40" "  +4.43%  "

# Row 41
Set-TextValue "B41" "InjectiveProtocol"
Set-TextValue "C41" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D41" "36.22"
Set-TextValue "E41" "  +2.35%  "

# Row 42
Set-TextValue "D42" "3.152.06"
Set-TextValue "E42" "  +2.69%  "

# Row 43
Set-TextValue "E43" "  +3.25%  "

# Row 44
Set-TextValue "D44" "0.0428"
Set-TextValue "E44" "  +2.47%  "

# Row 45
Set-TextValue "E45" "  +3.91%  "

# Row 46
Set-TextValue "B46" "dogwifhat"
Set-TextValue "C46" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D46" "2.80"
Set-TextValue "E46" "  +22.00%  "

# Row 47
Set-TextValue "B47" "ApeXProtocol"
Set-TextValue "C47" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D47" "3.23"
Set-TextValue "E47" "  +2.98%  "

# Row 48
Set-TextValue "B48" "Stellar"
Set-TextValue "C48" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D48" "0.135"
Set-TextValue "E48" "  +0.74%  "

# Row 49
Set-TextValue "B49" "THORChain"
Set-TextValue "C49" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D49" "8.73"
Set-TextValue "E49" "  +6.71%  "

# Row 50
Set-TextValue "B50" "FirstDigitalUSD"
Set-TextValue "C50" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D50" "0.999"
Set-TextValue "E50" "  -0.13%  "

# Row 51
Set-TextValue "D51" "140.81"
Set-TextValue "E51" "  +0.96%  "
